$d = $word.ActiveDocument

# --- Edit 1: "Create a webpage" bullet ---
# Split " with a list layout to display tasks." into a trailing space run
# plus a new run with the replacement sentence.
$find = $d.Content.Find
$find.Execute(" with a list layout to display tasks.", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null

$pWebpage = $d.Paragraphs(4)
$insertPos = $pWebpage.Range.End - 1
$ir = $d.Range($insertPos, $insertPos)
$ir.InsertAfter("using HTML to access API on web browser with interactions (buttons).")

# --- Edit 2: "Use JavaScript ..." bullet -> single replacement sentence ---
$pJs = $d.Paragraphs(5)
$insertPos2 = $pJs.Range.End - 1
$ir2 = $d.Range($insertPos2, $insertPos2)
$ir2.InsertAfter("Reformat so responses are presented in better format instead of JSON, perhaps CSS or HTML.")
$oldRange2 = $d.Range($pJs.Range.Start, $insertPos2)
$oldRange2.Delete()

# --- Edit 3: "Render tasks dynamically ..." bullet loses its text (stays as an
# empty bulleted paragraph); "Add interactivity ..." bullet and the trailing
# empty paragraph are removed entirely. ---
$pRender = $d.Paragraphs(6)
$renderText = $d.Range($pRender.Range.Start, $pRender.Range.End - 1)
$renderText.Delete()

$pAdd = $d.Paragraphs(7)
$pTrailing = $d.Paragraphs(8)
$removeRange = $d.Range($pAdd.Range.Start, $pTrailing.Range.End)
$removeRange.Delete()
